$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.330172223532543
$ws.Range("C2").Value = 0.1454990197619281
$ws.Range("D2").Value = 0.1627395962499705
$ws.Range("F2").Value = 1.449243925161589
$ws.Range("G2").Value = 0.002460514840477119
$ws.Range("I2").Value = 0.8401470047134261
$ws.Range("J2").Value = 0.1661648848737372
$ws.Range("L2").Value = 0.4100492350188745
$ws.Range("N2").Value = 1.31285967906598
$ws.Range("O2").Value = 3.5491773116849

$ws.Range("B3").Value = 1.234824730257856
$ws.Range("C3").Value = 0.1325167916962755
$ws.Range("D3").Value = 0.1621084604897121
$ws.Range("F3").Value = 1.451221181541058
$ws.Range("G3").Value = 0.002463331043428243
$ws.Range("I3").Value = 0.8473107237923934
$ws.Range("J3").Value = 0.1675507071515465
$ws.Range("L3").Value = 0.4021374256020636
$ws.Range("N3").Value = 1.321882594484997
$ws.Range("O3").Value = 3.556666444458557

$ws.Range("B4").Value = 1.176518825988296
$ws.Range("C4").Value = 0.1244903002525746
$ws.Range("D4").Value = 0.1617684163544482
$ws.Range("F4").Value = 1.453220094237835
$ws.Range("G4").Value = 0.002465154225579413
$ws.Range("I4").Value = 0.8521635876200619
$ws.Range("J4").Value = 0.1684575823552752
$ws.Range("L4").Value = 0.3974353838752194
$ws.Range("N4").Value = 1.327908392823709
$ws.Range("O4").Value = 3.563325983067955

$ws.Range("B5").Value = 1.152820232276895
$ws.Range("C5").Value = 0.1212056876374561
$ws.Range("D5").Value = 0.1616418336731158
$ws.Range("F5").Value = 1.454232090026494
$ws.Range("G5").Value = 0.002465920899962988
$ws.Range("I5").Value = 0.8542553733955742
$ws.Range("J5").Value = 0.1688412304223661
$ws.Range("L5").Value = 0.3955586415738992
$ws.Range("N5").Value = 1.330486233920539
$ws.Range("O5").Value = 3.56655801421536

$ws.Range("B6").Value = 1.148888866378826
$ws.Range("C6").Value = 0.1206594545192416
$ws.Range("D6").Value = 0.1616215402333054
$ws.Range("F6").Value = 1.454412056808138
$ws.Range("G6").Value = 0.002466049640137923
$ws.Range("I6").Value = 0.8546096094480689
$ws.Range("J6").Value = 0.1689057862418046
$ws.Range("L6").Value = 0.3952493938781174
$ws.Range("N6").Value = 1.330921672282486
$ws.Range("O6").Value = 3.567125986993744

$ws.Range("B7").Value = 1.176198967018536
$ws.Range("C7").Value = 0.1244460582801281
$ws.Range("D7").Value = 0.1617666606098922
$ws.Range("F7").Value = 1.453232942930697
$ws.Range("G7").Value = 0.002465164469199365
$ws.Range("I7").Value = 0.8521913358140303
$ws.Range("J7").Value = 0.1684626993028271
$ws.Range("L7").Value = 0.3974099137790148
$ws.Range("N7").Value = 1.327942663187102
$ws.Range("O7").Value = 3.563367473255312

$ws.Range("B8").Value = 1.297248048980464
$ws.Range("C8").Value = 0.1410343315836542
$ws.Range("D8").Value = 0.1625121564931078
$ws.Range("F8").Value = 1.449762816580382
$ws.Range("G8").Value = 0.002461466395356253
$ws.Range("I8").Value = 0.8425227409265048
$ws.Range("J8").Value = 0.166631107191197
$ws.Range("L8").Value = 0.4072890091514267
$ws.Range("N8").Value = 1.315870062571456
$ws.Range("O8").Value = 3.551331706851556

$ws.Range("B9").Value = 1.536448333789338
$ws.Range("C9").Value = 0.1731188463359388
$ws.Range("D9").Value = 0.1643487647506774
$ws.Range("F9").Value = 1.44918449031168
$ws.Range("G9").Value = 0.002454957301400221
$ws.Range("I9").Value = 0.8271692400991775
$ws.Range("J9").Value = 0.1634828561049595
$ws.Range("L9").Value = 0.4278913796481589
$ws.Range("N9").Value = 1.296043243790862
$ws.Range("O9").Value = 3.54409208856373

$ws.Range("B10").Value = 1.713232598583033
$ws.Range("C10").Value = 0.1964143019377218
$ws.Range("D10").Value = 0.1659241216900469
$ws.Range("F10").Value = 1.452556611079743
$ws.Range("G10").Value = 0.002450623423602092
$ws.Range("I10").Value = 0.818090525418036
$ws.Range("J10").Value = 0.1614392933172297
$ws.Range("L10").Value = 0.4437699733989149
$ws.Range("N10").Value = 1.283814056159791
$ws.Range("O10").Value = 3.548765028561007

$ws.Range("B11").Value = 1.793869737531736
$ws.Range("C11").Value = 0.2069507660854413
$ws.Range("D11").Value = 0.1666893714360924
$ws.Range("F11").Value = 1.454915578022636
$ws.Range("G11").Value = 0.002448748232081838
$ws.Range("I11").Value = 0.8144390627631424
$ws.Range("J11").Value = 0.1605679409897238
$ws.Range("L11").Value = 0.4511532148845561
$ws.Range("N11").Value = 1.278756691723103
$ws.Range("O11").Value = 3.553064548173182

$ws.Range("B12").Value = 1.824434629176551
$ws.Range("C12").Value = 0.2109317733121543
$ws.Range("D12").Value = 0.1669860954404427
$ws.Range("F12").Value = 1.455927477015777
$ws.Range("G12").Value = 0.002448051922956835
$ws.Range("I12").Value = 0.8131252086087812
$ws.Range("J12").Value = 0.160246348798907
$ws.Range("L12").Value = 0.4539718976102876
$ws.Range("N12").Value = 1.276914209311926
$ws.Range("O12").Value = 3.555005475567469

$ws.Range("B13").Value = 1.81785065679469
$ws.Range("C13").Value = 0.210074792000114
$ws.Range("D13").Value = 0.1669218827265269
$ws.Range("F13").Value = 1.455704271610514
$ws.Range("G13").Value = 0.002448201273414637
$ws.Range("I13").Value = 0.8134051062503787
$ws.Range("J13").Value = 0.1603152373270191
$ws.Range("L13").Value = 0.4533638325312381
$ws.Range("N13").Value = 1.277307792013943
$ws.Range("O13").Value = 3.554573547319535

$ws.Range("B14").Value = 1.796383751843223
$ws.Range("C14").Value = 0.207278465845377
$ws.Range("D14").Value = 0.1667136442841155
$ws.Range("F14").Value = 1.454996450572537
$ws.Range("G14").Value = 0.0024486906703006
$ws.Range("I14").Value = 0.8143295904271497
$ws.Range("J14").Value = 0.160541315714724
$ws.Range("L14").Value = 0.4513846536012807
$ws.Range("N14").Value = 1.278603654466451
$ws.Range("O14").Value = 3.553217959431862

$ws.Range("B15").Value = 1.783238420744169
$ws.Range("C15").Value = 0.2055644677543285
$ws.Range("D15").Value = 0.1665869946117056
$ws.Range("F15").Value = 1.45457833557461
$ws.Range("G15").Value = 0.002448992233392291
$ws.Range("I15").Value = 0.8149048353566641
$ws.Range("J15").Value = 0.160680884990315
$ws.Range("L15").Value = 0.4501753139778799
$ws.Range("N15").Value = 1.27940686368958
$ws.Range("O15").Value = 3.552428363639109

$ws.Range("B16").Value = 1.707966981225866
$ws.Range("C16").Value = 0.1957244807205996
$ws.Range("D16").Value = 0.1658750847560242
$ws.Range("F16").Value = 1.45241904584428
$ws.Range("G16").Value = 0.002450747904238741
$ws.Range("I16").Value = 0.8183387914515379
$ws.Range("J16").Value = 0.1614974100360218
$ws.Range("L16").Value = 0.4432906637674137
$ws.Range("N16").Value = 1.284154735651768
$ws.Range("O16").Value = 3.548527804439715

$ws.Range("B17").Value = 1.661844668426568
$ws.Range("C17").Value = 0.1896722736602499
$ws.Range("D17").Value = 0.1654507663300819
$ws.Range("F17").Value = 1.451305672233204
$ws.Range("G17").Value = 0.00245184957343578
$ws.Range("I17").Value = 0.820568014349071
$ws.Range("J17").Value = 0.1620132402328469
$ws.Range("L17").Value = 0.4391079908652102
$ws.Range("N17").Value = 1.287196867968163
$ws.Range("O17").Value = 3.546691841121543

$ws.Range("B18").Value = 1.635336878495764
$ws.Range("C18").Value = 0.1861854967749821
$ws.Range("D18").Value = 0.1652112881536567
$ws.Range("F18").Value = 1.450742940734415
$ws.Range("G18").Value = 0.002452492293311792
$ws.Range("I18").Value = 0.8218952380197138
$ws.Range("J18").Value = 0.1623154177824144
$ws.Range("L18").Value = 0.4367172987732886
$ws.Range("N18").Value = 1.288994231388614
$ws.Range("O18").Value = 3.545840437384982

$ws.Range("B19").Value = 1.626365381573748
$ws.Range("C19").Value = 0.1850039594608859
$ws.Range("D19").Value = 0.1651309928497966
$ws.Range("F19").Value = 1.450565746869756
$ws.Range("G19").Value = 0.002452711466785575
$ws.Range("I19").Value = 0.8223523456109376
$ws.Range("J19").Value = 0.1624186723767824
$ws.Range("L19").Value = 0.4359104459233691
$ws.Range("N19").Value = 1.289610967693477
$ws.Range("O19").Value = 3.545587299590494

$ws.Range("B20").Value = 1.666752354566313
$ws.Range("C20").Value = 0.1903171334856779
$ws.Range("D20").Value = 0.1654954622486358
$ws.Range("F20").Value = 1.451416156256343
$ws.Range("G20").Value = 0.00245173136070459
$ws.Range("I20").Value = 0.820326048452678
$ws.Range("J20").Value = 0.1619577615536301
$ws.Range("L20").Value = 0.4395516851993904
$ws.Range("N20").Value = 1.286868101100453
$ws.Range("O20").Value = 3.546866106354997

$ws.Range("B21").Value = 1.802688320252912
$ws.Range("C21").Value = 0.208100058275221
$ws.Range("D21").Value = 0.1667746209951773
$ws.Range("F21").Value = 1.455201135910798
$ws.Range("G21").Value = 0.002448546548940128
$ws.Range("I21").Value = 0.8140561772016994
$ws.Range("J21").Value = 0.1604746839294098
$ws.Range("L21").Value = 0.4519653689564649
$ws.Range("N21").Value = 1.278221057782055
$ws.Range("O21").Value = 3.553607637715402

$ws.Range("B22").Value = 1.891700587544278
$ws.Range("C22").Value = 0.219670130761358
$ws.Range("D22").Value = 0.1676510568967515
$ws.Range("F22").Value = 1.458366175395099
$ws.Range("G22").Value = 0.002446545412640238
$ws.Range("I22").Value = 0.8103599451504095
$ws.Range("J22").Value = 0.1595541894055277
$ws.Range("L22").Value = 0.4602112883755467
$ws.Range("N22").Value = 1.272993027642691
$ws.Range("O22").Value = 3.559836907664646

$ws.Range("B23").Value = 1.844178100531508
$ws.Range("C23").Value = 0.2134997937420735
$ws.Range("D23").Value = 0.1671796034807116
$ws.Range("F23").Value = 1.456613683864205
$ws.Range("G23").Value = 0.002447606128267004
$ws.Range("I23").Value = 0.8122959344240428
$ws.Range("J23").Value = 0.1600410142704494
$ws.Range("L23").Value = 0.455798194346599
$ws.Range("N23").Value = 1.275744622813505
$ws.Range("O23").Value = 3.556345330178033

$ws.Range("B24").Value = 1.664533562910833
$ws.Range("C24").Value = 0.19002561503018
$ws.Range("D24").Value = 0.1654752413027865
$ws.Range("F24").Value = 1.451365965449739
$ws.Range("G24").Value = 0.002451784775495818
$ws.Range("I24").Value = 0.8204352991543047
$ws.Range("J24").Value = 0.1619828259581055
$ws.Range("L24").Value = 0.4393510474241964
$ws.Range("N24").Value = 1.287016585826123
$ws.Range("O24").Value = 3.54678668524457

$ws.Range("B25").Value = 1.471550133146195
$ws.Range("C25").Value = 0.1644873729405276
$ws.Range("D25").Value = 0.1638120509611767
$ws.Range("F25").Value = 1.448674227336141
$ws.Range("G25").Value = 0.002456639125739827
$ws.Range("I25").Value = 0.830936367040124
$ws.Range("J25").Value = 0.1642871566283315
$ws.Range("L25").Value = 0.4221870430141053
$ws.Range("N25").Value = 1.300995863632558
$ws.Range("O25").Value = 3.544297133449135
